# Update financial time series calibration estimates (std moving average)
# across the AR, SETAR, GARCH, TARCH and AR_TARCH worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.007093537996962009
$ws.Range("B3").Value = 0.7713960392002097
$ws.Range("B4").Value = 0.1044958812896311
$ws.Range("B5").Value = "[1.0, 0.020742072449138466, -0.046057578182718904, -0.09987740779527826, -0.11362098275977535, 0.07053872377467447, 0.027242553341174296, 0.038508958376020745, 0.022524839067012383, -0.014046769068594755, -0.041041619122136656, -0.030166978710828258, -0.0071719504665959215, 0.05899539513191117, 0.05611133351917008, 0.03036155380748464, -0.006839678743164695, -0.03382079672006651, -0.04097821460225061, 0.002374139634956119]"

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1855753646724737
$ws.Range("B5").Value = 0.5381593286732445
$ws.Range("B6").Value = 0.06376423619064685
$ws.Range("B7").Value = 0.1782192630903307
$ws.Range("B8").Value = 0.5728745678554105
$ws.Range("B9").Value = 0.06179907510183994
$ws.Range("B10").Value = "[1.0, 0.059786836446393266, 0.008164407218926334, -0.0022021570605025636, -0.00600966558310538, 0.010383753667678432, 0.016637941679735126, -0.02002175112349545, 0.00828413956072075, -0.008164560034633245, 0.012766997876978145, 0.018733138120191988, 0.013428564824778264, 0.029600391644876812, 0.019874002882877304, -4.271655568156782e-05, -0.015510411516043895, -0.007215655497153533, -0.004023081281778683, -0.012554653841117053]"

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.001042451192802102
$ws.Range("B3").Value = 0.1127127736981839
$ws.Range("B4").Value = 0.04456470552019402
$ws.Range("B5").Value = 0.000000008487559062106518
$ws.Range("B6").Value = "[1.0, 0.0046055459593722405, -0.03433958921631113, -0.09108536610602833, -0.13236200747524696, 0.09214281163745758, 0.05071434479689837, 0.04263232309750947, 0.03213706187886878, 8.049679111961297e-05, -0.053645870199857094, -0.025443903188358093, -0.010255174641980542, 0.05561505926299536, 0.04912782792613663, 0.018314890859462658, -0.0050686743117231165, -0.0302937502907785, -0.05743611288880337, -0.0029958877751506195]"

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.001327744496415256
$ws.Range("B3").Value = 0.1122934717130698
$ws.Range("B4").Value = 0.008452274866491095
$ws.Range("B5").Value = 0.000000005196964653699701
$ws.Range("B6").Value = "[1.0, 0.0029458747712240037, -0.034465194544982605, -0.09075287544891736, -0.13238801623191018, 0.09284870902027942, 0.050533654351551584, 0.041719372478274556, 0.03344411420821392, 0.00010641221794151788, -0.05233931341253185, -0.025471747071868847, -0.01021390984255276, 0.05487392018786231, 0.04860478141715638, 0.018248576979006625, -0.004905016845498029, -0.030547681261132508, -0.05714795099327387, -0.0032193038803278273]"
$ws.Range("B7").Value = 0.08014793525174554

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.007170031704656806
$ws.Range("B3").Value = 0.1017035991003381
$ws.Range("B4").Value = 0.000000000016376280646107
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, -0.002525254187664789, -0.04395351866747207, -0.09883537232329236, -0.11335429257458329, 0.07431803183308167, 0.025486775210797865, 0.03777172350984037, 0.021905212588827754, -0.012951354087596116, -0.03879417246439179, -0.029310032592469414, -0.0077053457037148335, 0.05671236126776432, 0.054249197375603596, 0.029100293210020035, -0.006396670870171413, -0.03277595465232519, -0.04002549926969663, 0.0029365595332461053]"
$ws.Range("B7").Value = 0.05286187043770072
$ws.Range("B9").Value = 0.7727068425809218
